# Applies the cryptos-list price/volume refresh described in the commit.
# All target cells are plain text in the sheet (t="inlineStr"), so the
# NumberFormat is forced to "@" (Text) before writing each value -- this
# stops Excel/COM from auto-coercing numeric-looking strings (e.g. "1.00",
# "0.0000140") into actual numbers and losing the formatting / precision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.296.81'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.08%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.573.21'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.15%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '505.78'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.54'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.56%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.577'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -5.70%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.577.17'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.56'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +7.45%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.95%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.91%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.035.07'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.463.88'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.70%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.53'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.22%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000140'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.599.47'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.75%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.72%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '344.92'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.41'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.12'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.996'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '59.92'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.419'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.71%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.166'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.64%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.94%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0839'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.35'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '19.31'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '154.13'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.55'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.99%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.70'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.32%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.96'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.34%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.60%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.849'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +8.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.848'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.54%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.74'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '35.90'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '293.18'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.84%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.620'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.18%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0992'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.22%  '
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0554'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.76%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.78'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.98%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.86'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0232'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.39%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.995.96'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.41%  '
